$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.973.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.786.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.568.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.970.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0707"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.81%  "

$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  +1.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("E35").Value = "  +2.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.33%  "

$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("E39").Value = "  +1.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.13%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0959"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
